$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.022893260322361
$ws.Range("D2").Value = 1.03435723317423
$ws.Range("E2").Value = 1.02358069974562
$ws.Range("F2").Value = 1.045237848330065
$ws.Range("I2").Value = 1.032963650995014
$ws.Range("J2").Value = 1.028076539021879
$ws.Range("K2").Value = 1.03715693959729
$ws.Range("L2").Value = 1.026411718795603
$ws.Range("M2").Value = 1.048006650703347
$ws.Range("N2").Value = 1.013406421555046
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.023743724085129
$ws.Range("D3").Value = 1.035016956119125
$ws.Range("E3").Value = 1.024299909772917
$ws.Range("F3").Value = 1.04605657031761
$ws.Range("I3").Value = 1.03309378783471
$ws.Range("J3").Value = 1.028565761338694
$ws.Range("K3").Value = 1.037626231876103
$ws.Range("L3").Value = 1.026938077090544
$ws.Range("M3").Value = 1.048636753558684
$ws.Range("N3").Value = 1.013570117490231
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.024294745204023
$ws.Range("D4").Value = 1.035444271663311
$ws.Range("E4").Value = 1.024766293235082
$ws.Range("F4").Value = 1.046587101187107
$ws.Range("I4").Value = 1.033176842504439
$ws.Range("J4").Value = 1.02888239956211
$ws.Range("K4").Value = 1.037929642554075
$ws.Range("L4").Value = 1.027279002768743
$ws.Range("M4").Value = 1.049044583764239
$ws.Range("N4").Value = 1.013676015285077
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.0245265638726
$ws.Range("D5").Value = 1.035624016373543
$ws.Range("E5").Value = 1.024962599992224
$ws.Range("F5").Value = 1.046810317005303
$ws.Range("I5").Value = 1.033211482061958
$ws.Range("J5").Value = 1.029015531833404
$ws.Range("K5").Value = 1.038057134479515
$ws.Range("L5").Value = 1.027422407150525
$ws.Range("M5").Value = 1.049216060566596
$ws.Range("N5").Value = 1.013720528398011
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.024565497134087
$ws.Range("D6").Value = 1.035654202185109
$ws.Range("E6").Value = 1.0249955747319
$ws.Range("F6").Value = 1.046847806457421
$ws.Range("I6").Value = 1.033217281954336
$ws.Range("J6").Value = 1.02903788632044
$ws.Range("K6").Value = 1.038078537262816
$ws.Range("L6").Value = 1.027446489983236
$ws.Range("M6").Value = 1.04924485369144
$ws.Range("N6").Value = 1.013728001961157
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.024297842112817
$ws.Range("D7").Value = 1.035446673025392
$ws.Range("E7").Value = 1.024768915358154
$ws.Range("F7").Value = 1.046590083100543
$ws.Range("I7").Value = 1.033177306447724
$ws.Range("J7").Value = 1.028884178414624
$ws.Range("K7").Value = 1.037931346352501
$ws.Range("L7").Value = 1.027280918635528
$ws.Range("M7").Value = 1.04904687494965
$ws.Range("N7").Value = 1.013676610096899
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.023180529697011
$ws.Range("D8").Value = 1.034580099343518
$ws.Range("E8").Value = 1.023823550658387
$ws.Range("F8").Value = 1.045514380059126
$ws.Range("I8").Value = 1.033007869575572
$ws.Range("J8").Value = 1.028241856940187
$ws.Range("K8").Value = 1.037315590623056
$ws.Range("L8").Value = 1.026589533512522
$ws.Range("M8").Value = 1.048219572789953
$ws.Range("N8").Value = 1.013461748134709
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.021217222736984
$ws.Range("D9").Value = 1.033056466122527
$ws.Range("E9").Value = 1.02216548547548
$ws.Range("F9").Value = 1.043624786914664
$ws.Range("I9").Value = 1.032700503781486
$ws.Range("J9").Value = 1.027110663190856
$ws.Range("K9").Value = 1.036228676080106
$ws.Range("L9").Value = 1.025373863668907
$ws.Range("M9").Value = 1.046762690641952
$ws.Range("N9").Value = 1.013082967789169
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.019912171071562
$ws.Range("D10").Value = 1.032043093792246
$ws.Range("E10").Value = 1.021065452718867
$ws.Range("F10").Value = 1.042369167825394
$ws.Range("I10").Value = 1.032489720758315
$ws.Range("J10").Value = 1.026357052738459
$ws.Range("K10").Value = 1.03550288975453
$ws.Range("L10").Value = 1.02456527326754
$ws.Range("M10").Value = 1.045792163400002
$ws.Range("N10").Value = 1.012830364083205
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.019347995151983
$ws.Range("D11").Value = 1.031604880971045
$ws.Range("E11").Value = 1.020590416698637
$ws.Range("F11").Value = 1.041826471129601
$ws.Range("I11").Value = 1.032397065285407
$ws.Range("J11").Value = 1.02603087010826
$ws.Range("K11").Value = 1.035188354804632
$ws.Range("L11").Value = 1.024215602223671
$ws.Range("M11").Value = 1.045372108533324
$ws.Range("N11").Value = 1.012720970041442
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.019138574918416
$ws.Range("D12").Value = 1.031442198859117
$ws.Range("E12").Value = 1.020414162141512
$ws.Range("F12").Value = 1.041625040909979
$ws.Range("I12").Value = 1.032362441523011
$ws.Range("J12").Value = 1.025909733024548
$ws.Range("K12").Value = 1.035071484170041
$ws.Range("L12").Value = 1.024085788423045
$ws.Range("M12").Value = 1.045216111712624
$ws.Range("N12").Value = 1.012680334453323
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.019183489916502
$ws.Range("D13").Value = 1.0314770906376
$ws.Range("E13").Value = 1.020451960489736
$ws.Range("F13").Value = 1.04166824148723
$ws.Range("I13").Value = 1.032369877811497
$ws.Range("J13").Value = 1.02593571633903
$ws.Range("K13").Value = 1.035096555030558
$ws.Range("L13").Value = 1.024113630752264
$ws.Range("M13").Value = 1.045249572162971
$ws.Range("N13").Value = 1.012689050997339
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.019330681548571
$ws.Range("D14").Value = 1.031591431777353
$ws.Range("E14").Value = 1.020575843453898
$ws.Range("F14").Value = 1.041809817742408
$ws.Range("I14").Value = 1.032394207501185
$ws.Range("J14").Value = 1.026020856429667
$ws.Range("K14").Value = 1.035178695015652
$ws.Range("L14").Value = 1.02420487034665
$ws.Range("M14").Value = 1.045359213160116
$ws.Range("N14").Value = 1.012717611124396
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.019421389746656
$ws.Range("D15").Value = 1.031661893062091
$ws.Range("E15").Value = 1.020652197727768
$ws.Range("F15").Value = 1.041897067682249
$ws.Range("I15").Value = 1.032409170364925
$ws.Range("J15").Value = 1.026073316959218
$ws.Range("K15").Value = 1.035229299132083
$ws.Range("L15").Value = 1.024261095335663
$ws.Range("M15").Value = 1.045426770658788
$ws.Range("N15").Value = 1.012735207741021
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.019949633089913
$ws.Range("D16").Value = 1.032072189001282
$ws.Range("E16").Value = 1.021097006535509
$ws.Range("F16").Value = 1.042405205949965
$ws.Range("I16").Value = 1.032495840874828
$ws.Range("J16").Value = 1.026378703373505
$ws.Range("K16").Value = 1.035523758947597
$ws.Range("L16").Value = 1.024588489466149
$ws.Range("M16").Value = 1.045820045206151
$ws.Range("N16").Value = 1.012837623930958
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.020281233634186
$ws.Range("D17").Value = 1.032329714514602
$ws.Range("E17").Value = 1.021376368800764
$ws.Range("F17").Value = 1.042724215413277
$ws.Range("I17").Value = 1.032549836670903
$ws.Range("J17").Value = 1.026570301433059
$ws.Range("K17").Value = 1.035708395947325
$ws.Range("L17").Value = 1.024793977671454
$ws.Range("M17").Value = 1.0460667879915
$ws.Range("N17").Value = 1.012901863231361
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.020474739287133
$ws.Range("D18").Value = 1.032479981018716
$ws.Range("E18").Value = 1.021539439984225
$ws.Range("F18").Value = 1.042910384067697
$ws.Range("I18").Value = 1.032581197823208
$ws.Range("J18").Value = 1.026682070343464
$ws.Range("K18").Value = 1.035816065952611
$ws.Range("L18").Value = 1.02491387916971
$ws.Range("M18").Value = 1.04621072708321
$ws.Range("N18").Value = 1.012939331451952
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.020540734704966
$ws.Range("D19").Value = 1.032531227493321
$ws.Range("E19").Value = 1.021595064002974
$ws.Range("F19").Value = 1.042973878959827
$ws.Range("I19").Value = 1.032591868469896
$ws.Range("J19").Value = 1.026720182820107
$ws.Range("K19").Value = 1.035852774251113
$ws.Range("L19").Value = 1.024954769852691
$ws.Range("M19").Value = 1.046259809627802
$ws.Range("N19").Value = 1.012952106881895
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.020245646851409
$ws.Range("D20").Value = 1.03230207863454
$ws.Range("E20").Value = 1.021346383044584
$ws.Range("F20").Value = 1.042689978789747
$ws.Range("I20").Value = 1.032544057255367
$ws.Range("J20").Value = 1.026549743420207
$ws.Range("K20").Value = 1.035688588785749
$ws.Range("L20").Value = 1.024771926201095
$ws.Range("M20").Value = 1.046040312930762
$ws.Range("N20").Value = 1.012894971111646
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.019287333408968
$ws.Range("D21").Value = 1.031557758671877
$ws.Range("E21").Value = 1.02053935761162
$ws.Range("F21").Value = 1.041768122865354
$ws.Range("I21").Value = 1.032387048733514
$ws.Range("J21").Value = 1.025995784186557
$ws.Range("K21").Value = 1.035154507884595
$ws.Range("L21").Value = 1.024178000627849
$ws.Range("M21").Value = 1.045326925769905
$ws.Range("N21").Value = 1.012709200923844
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.018685613476223
$ws.Range("D22").Value = 1.031090294739516
$ws.Range("E22").Value = 1.020033077652655
$ws.Range("F22").Value = 1.041189392990374
$ws.Range("I22").Value = 1.03228713166781
$ws.Range("J22").Value = 1.025647614032941
$ws.Range("K22").Value = 1.034818488846245
$ws.Range("L22").Value = 1.023804979387594
$ws.Range("M22").Value = 1.044878566881527
$ws.Range("N22").Value = 1.01259238982822
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.019004519308053
$ws.Range("D23").Value = 1.031338056335392
$ws.Range("E23").Value = 1.020301358607829
$ws.Range("F23").Value = 1.041496104779108
$ws.Range("I23").Value = 1.032340213069714
$ws.Range("J23").Value = 1.025832173246045
$ws.Range("K23").Value = 1.034996639388703
$ws.Range("L23").Value = 1.024002686307115
$ws.Range("M23").Value = 1.045116233154382
$ws.Range("N23").Value = 1.012654314388767
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.020261726712955
$ws.Range("D24").Value = 1.032314565925669
$ws.Range("E24").Value = 1.021359931934487
$ws.Range("F24").Value = 1.042705448550376
$ws.Range("I24").Value = 1.032546669137613
$ws.Range("J24").Value = 1.026559032661284
$ws.Range("K24").Value = 1.035697538868851
$ws.Range("L24").Value = 1.024781890176491
$ws.Range("M24").Value = 1.046052275815084
$ws.Range("N24").Value = 1.012898085368499
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.021724118353664
$ws.Range("D25").Value = 1.03344994986019
$ws.Range("E25").Value = 1.022593201023993
$ws.Range("F25").Value = 1.044112576506241
$ws.Range("I25").Value = 1.032781003330483
$ws.Range("J25").Value = 1.027403018012414
$ws.Range("K25").Value = 1.036509883045092
$ws.Range("L25").Value = 1.02568782228751
$ws.Range("M25").Value = 1.047139208853089
$ws.Range("N25").Value = 1.013180908364897
